$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '48.134.47'
$ws.Cells.Item(2, 5).Value = '  +0.70%  '
$ws.Cells.Item(3, 4).Value = '2.503.57'
$ws.Cells.Item(3, 5).Value = '  +0.35%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '321.23'
$ws.Cells.Item(5, 5).Value = '  -0.58%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '107.63'
$ws.Cells.Item(6, 5).Value = '  -1.20%  '
$ws.Cells.Item(8, 5).Value = '  -0.06%  '
$ws.Cells.Item(9, 5).Value = '  -2.03%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '39.65'
$ws.Cells.Item(10, 5).Value = '  -2.53%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '20.14'
$ws.Cells.Item(11, 5).Value = '  +8.00%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.0812'
$ws.Cells.Item(12, 5).Value = '  -0.13%  '
$ws.Cells.Item(13, 5).Value = '  -0.13%  '
$ws.Cells.Item(14, 5).Value = '  -1.43%  '
$ws.Cells.Item(15, 4).Value = '2.895.97'
$ws.Cells.Item(15, 5).Value = '  +0.46%  '
$ws.Cells.Item(16, 4).Value = '2.507.67'
$ws.Cells.Item(16, 5).Value = '  +0.08%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '0.834'
$ws.Cells.Item(17, 5).Value = '  -1.98%  '
$ws.Cells.Item(18, 4).Value = '48.004.78'
$ws.Cells.Item(18, 5).Value = '  +0.67%  '
$ws.Cells.Item(19, 5).Value = '  -1.84%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '6.74'
$ws.Cells.Item(20, 5).Value = '  +1.54%  '
$ws.Cells.Item(21, 5).Value = '  +0.12%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '276.99'
$ws.Cells.Item(23, 5).Value = '  +12.06%  '
$ws.Cells.Item(24, 5).Value = '  +1.07%  '
$ws.Cells.Item(25, 5).Value = '  -0.65%  '
$ws.Cells.Item(26, 5).Value = '  -0.07%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '25.96'
$ws.Cells.Item(27, 5).Value = '  +0.41%  '
$ws.Cells.Item(28, 5).Value = '  -2.72%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.141'
$ws.Cells.Item(29, 5).Value = '  +1.08%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '35.26'
$ws.Cells.Item(30, 5).Value = '  -0.10%  '
$ws.Cells.Item(31, 5).Value = '  -4.60%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '49.68'
$ws.Cells.Item(32, 5).Value = '  -0.24%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '19.56'
$ws.Cells.Item(33, 5).Value = '  -1.71%  '
$ws.Cells.Item(34, 5).Value = '  -0.12%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.0784'
$ws.Cells.Item(36, 5).Value = '  -0.81%  '
$ws.Cells.Item(37, 5).Value = '  -0.96%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '4.63'
$ws.Cells.Item(38, 5).Value = '  -0.83%  '
$ws.Cells.Item(39, 5).Value = '  -2.65%  '
$ws.Cells.Item(40, 5).Value = '  -0.49%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '121.24'
$ws.Cells.Item(41, 5).Value = '  +1.82%  '
$ws.Cells.Item(42, 5).Value = '  -0.61%  '
$ws.Cells.Item(43, 5).Value = '  -5.12%  '
$ws.Cells.Item(44, 5).Value = '  +1.43%  '
$ws.Cells.Item(45, 4).Value = '2.024.56'
$ws.Cells.Item(45, 5).Value = '  +1.07%  '
$ws.Cells.Item(46, 5).Value = '  +3.15%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.00'
$ws.Cells.Item(47, 5).Value = '  -1.59%  '
$ws.Cells.Item(48, 5).Value = '  +1.47%  '
$ws.Cells.Item(49, 5).Value = '  -0.30%  '
$ws.Cells.Item(50, 5).Value = '  +0.72%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '80.37'
$ws.Cells.Item(51, 5).Value = '  +3.33%  '
